$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------------

function Wrap-Body($bodyXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">' + `
        '<w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

function Get-ParaByText($text) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $rng.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    return $rng.Paragraphs(1)
}

# ---------------------------------------------------------------------------
# 1) Topic title paragraph: replace text + formatting, then add a blank
#    paragraph right after it (no pPr / no run at all).
# ---------------------------------------------------------------------------

$pTopic = Get-ParaByText("Application of IOT in Logistics & Supply Chain Management")
$titleXml = '<w:p><w:r><w:t>Proposal of protocol based vehicle</w:t></w:r>' + `
            '<w:r><w:t xml:space="preserve"> tracking system</w:t></w:r>' + `
            '<w:r><w:t xml:space="preserve"> using ESP32.</w:t></w:r></w:p>' + `
            '<w:p></w:p>' + `
            '<w:p><w:r><w:t>ZZTOPICDUMMYZZ</w:t></w:r></w:p>'
$pTopic.Range.InsertXML((Wrap-Body $titleXml))

$pDummy1 = Get-ParaByText("ZZTOPICDUMMYZZ")
$pDummy1.Range.Delete()

# ---------------------------------------------------------------------------
# 2) Problem definition paragraph.
# ---------------------------------------------------------------------------

$pProblem = Get-ParaByText("The current logistic management system lacks a layer based framework which addresses different aspects thoroughly in the supply chain. So a proper management framework needs to be proposed in order to improve the quality of the service.")
$hl = '<w:highlight w:val="white"/>'
$problemXml = '<w:p><w:r><w:rPr>' + $hl + '</w:rPr><w:t xml:space="preserve"> The </w:t></w:r>' + `
              '<w:r><w:rPr>' + $hl + '</w:rPr><w:t>current vehicle tracking systems lacks the proper layer based design approach which addresses different aspects of the tracking. Our layer based model proposes new approach in the vehicle tracking systems.</w:t></w:r></w:p>'
$pProblem.Range.InsertXML((Wrap-Body $problemXml))

# ---------------------------------------------------------------------------
# 3) Objectives paragraph.
# ---------------------------------------------------------------------------

$pObjective = Get-ParaByText("The main objective of the project is to overcome the challenges faced with Logistics & Supply Chain.")
$objectiveXml = '<w:p><w:r><w:rPr>' + $hl + '</w:rPr><w:t>The main objective of the project is to overcome the challenges faced with Logistics &amp; Supply Chain</w:t></w:r>' + `
                '<w:r><w:rPr>' + $hl + '</w:rPr><w:t xml:space="preserve"> with the new layer based design</w:t></w:r>' + `
                '<w:r><w:rPr>' + $hl + '</w:rPr><w:t>.</w:t></w:r></w:p>'
$pObjective.Range.InsertXML((Wrap-Body $objectiveXml))

# ---------------------------------------------------------------------------
# 4) First project-life-cycle bullet.
# ---------------------------------------------------------------------------

$pBullet1 = Get-ParaByText("Finding the loopholes in the existing logistics management technology.")
$numPr = '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>'
$bullet1Xml = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/>' + $numPr + '<w:rPr>' + $hl + '</w:rPr></w:pPr>' + `
              '<w:r><w:rPr>' + $hl + '</w:rPr><w:t>Finding the loopholes in the existing</w:t></w:r>' + `
              '<w:r><w:rPr>' + $hl + '</w:rPr><w:t xml:space="preserve"> live tracking technology</w:t></w:r>' + `
              '<w:r><w:rPr>' + $hl + '</w:rPr><w:t>.</w:t></w:r></w:p>'
$pBullet1.Range.InsertXML((Wrap-Body $bullet1Xml))

# ---------------------------------------------------------------------------
# 5) Remove the "Proposing new methods in the logistics chain." bullet
#    entirely (it is deleted by the edit).
# ---------------------------------------------------------------------------

$pBullet2 = Get-ParaByText("Proposing new methods in the logistics chain.")
$pBullet2.Range.Delete()

# ---------------------------------------------------------------------------
# 6) After the "Proposing layer based approach..." bullet, insert a new
#    blank (non-list) paragraph that only carries the highlight rPr.
# ---------------------------------------------------------------------------

$pBullet3 = Get-ParaByText("Proposing layer based approach for the complete logistics management system in order to improve the quality of service.")
$bullet3Xml = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/>' + $numPr + '<w:rPr>' + $hl + '</w:rPr></w:pPr>' + `
              '<w:r><w:rPr>' + $hl + '</w:rPr><w:t>Proposing layer based approach for the complete logistics management system in order to improve the quality of service.</w:t></w:r></w:p>' + `
              '<w:p><w:pPr><w:rPr>' + $hl + '</w:rPr></w:pPr></w:p>' + `
              '<w:p><w:r><w:t>ZZBULLETDUMMYZZ</w:t></w:r></w:p>'
$pBullet3.Range.InsertXML((Wrap-Body $bullet3Xml))

$pDummy2 = Get-ParaByText("ZZBULLETDUMMYZZ")
$pDummy2.Range.Delete()

# ---------------------------------------------------------------------------
# 7) Append "Reference paper :" runs (+ proofErr tags) to the blank
#    paragraph that currently only holds a single space, then insert a new
#    paragraph after it containing the reference hyperlink.
# ---------------------------------------------------------------------------

$pRef = Get-ParaByText("Generic Project life cycle for chosen technology:")
# the blank-space paragraph we need sits right before the trailing
# " " paragraph at the very end of the body; locate it via the bookmark
# paragraph which is unique and immediately precedes it.
$pSpace = $d.Paragraphs($d.Paragraphs.Count - 1)

$spaceXml = '<w:p><w:pPr><w:rPr>' + $hl + '</w:rPr></w:pPr>' + `
            '<w:r><w:rPr>' + $hl + '</w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
            '<w:r><w:rPr>' + $hl + '</w:rPr><w:t xml:space="preserve">Reference </w:t></w:r>' + `
            '<w:proofErr w:type="gramStart"/>' + `
            '<w:r><w:rPr>' + $hl + '</w:rPr><w:t>paper :</w:t></w:r>' + `
            '<w:proofErr w:type="gramEnd"/></w:p>' + `
            '<w:p><w:pPr><w:rPr>' + $hl + '</w:rPr></w:pPr>' + `
            '<w:r><w:t>ZZHYPERLINKDUMMYZZ</w:t></w:r></w:p>'
$pSpace.Range.InsertXML((Wrap-Body $spaceXml))

# ---------------------------------------------------------------------------
# 8) Register the "Hyperlink" character style used by the new run, then
#    turn the placeholder paragraph into a real hyperlink field pointing at
#    the reference paper's URL.
# ---------------------------------------------------------------------------

$s = $d.Styles.Add("Hyperlink", 2)
$s.BaseStyle = $d.Styles("DefaultParagraphFont")
$s.Font.Color = 16711680
$s.Font.Underline = 1
$s.Priority = 99
$s.UnhideWhenUsed = $true

$url = "https://www.ijitee.org/wp-content/uploads/papers/v8i6/F3569048619.pdf"
$pDummy3 = Get-ParaByText("ZZHYPERLINKDUMMYZZ")
$linkRange = $pDummy3.Range
$linkRange.Text = $url
$linkRange = $pDummy3.Range
$linkRange.MoveEnd(1, -1) | Out-Null
$d.Hyperlinks.Add($linkRange, $url, "", "", $url) | Out-Null
